$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-16 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-17 Thursday", 2) | Out-Null
$d.Content.Find.Execute("46-35=", $true, $false, $false, $false, $false, $true, 1, $false, "5+91=", 2) | Out-Null
$d.Content.Find.Execute("45+37=", $true, $false, $false, $false, $false, $true, 1, $false, "77-67=", 2) | Out-Null
$d.Content.Find.Execute("36+32=", $true, $false, $false, $false, $false, $true, 1, $false, "60-31=", 2) | Out-Null
$d.Content.Find.Execute("61-13=", $true, $false, $false, $false, $false, $true, 1, $false, "89-53=", 2) | Out-Null
$d.Content.Find.Execute("45-22=", $true, $false, $false, $false, $false, $true, 1, $false, "32-9=", 2) | Out-Null
$d.Content.Find.Execute("61-40=", $true, $false, $false, $false, $false, $true, 1, $false, "72+22=", 2) | Out-Null
$d.Content.Find.Execute("10+7=", $true, $false, $false, $false, $false, $true, 1, $false, "27+58=", 2) | Out-Null
$d.Content.Find.Execute("54-4=", $true, $false, $false, $false, $false, $true, 1, $false, "45+47=", 2) | Out-Null
$d.Content.Find.Execute("41-0=", $true, $false, $false, $false, $false, $true, 1, $false, "74-39=", 2) | Out-Null
$d.Content.Find.Execute("63-26=", $true, $false, $false, $false, $false, $true, 1, $false, "38+42=", 2) | Out-Null
$d.Content.Find.Execute("56-16=", $true, $false, $false, $false, $false, $true, 1, $false, "32+0=", 2) | Out-Null
$d.Content.Find.Execute("57-17=", $true, $false, $false, $false, $false, $true, 1, $false, "37-15=", 2) | Out-Null
$d.Content.Find.Execute("63-39=", $true, $false, $false, $false, $false, $true, 1, $false, "67-53=", 2) | Out-Null
$d.Content.Find.Execute("71-68=", $true, $false, $false, $false, $false, $true, 1, $false, "21-19=", 2) | Out-Null
$d.Content.Find.Execute("20+57=", $true, $false, $false, $false, $false, $true, 1, $false, "42+49=", 2) | Out-Null
$d.Content.Find.Execute("60+22=", $true, $false, $false, $false, $false, $true, 1, $false, "81-74=", 2) | Out-Null
$d.Content.Find.Execute("88-78=", $true, $false, $false, $false, $false, $true, 1, $false, "65+20=", 2) | Out-Null
$d.Content.Find.Execute("25+54=", $true, $false, $false, $false, $false, $true, 1, $false, "95-76=", 2) | Out-Null
$d.Content.Find.Execute("92-0=", $true, $false, $false, $false, $false, $true, 1, $false, "87+6=", 2) | Out-Null
$d.Content.Find.Execute("92-23=", $true, $false, $false, $false, $false, $true, 1, $false, "96-89=", 2) | Out-Null
$d.Content.Find.Execute("16+54=", $true, $false, $false, $false, $false, $true, 1, $false, "73-70=", 2) | Out-Null
$d.Content.Find.Execute("24+8=", $true, $false, $false, $false, $false, $true, 1, $false, "8+28=", 2) | Out-Null
$d.Content.Find.Execute("31-13=", $true, $false, $false, $false, $false, $true, 1, $false, "63-11=", 2) | Out-Null
$d.Content.Find.Execute("44+6=", $true, $false, $false, $false, $false, $true, 1, $false, "58+21=", 2) | Out-Null
$d.Content.Find.Execute("18+41=", $true, $false, $false, $false, $false, $true, 1, $false, "16+3=", 2) | Out-Null
$d.Content.Find.Execute("93-87=", $true, $false, $false, $false, $false, $true, 1, $false, "70-18=", 2) | Out-Null
$d.Content.Find.Execute("67-24=", $true, $false, $false, $false, $false, $true, 1, $false, "32+4=", 2) | Out-Null
$d.Content.Find.Execute("75+6=", $true, $false, $false, $false, $false, $true, 1, $false, "96-55=", 2) | Out-Null
$d.Content.Find.Execute("59+28=", $true, $false, $false, $false, $false, $true, 1, $false, "40+30=", 2) | Out-Null
$d.Content.Find.Execute("78-51=", $true, $false, $false, $false, $false, $true, 1, $false, "8+65=", 2) | Out-Null
$d.Content.Find.Execute("57+10=", $true, $false, $false, $false, $false, $true, 1, $false, "99-73=", 2) | Out-Null
$d.Content.Find.Execute("57+41=", $true, $false, $false, $false, $false, $true, 1, $false, "43-40=", 2) | Out-Null
$d.Content.Find.Execute("49-43=", $true, $false, $false, $false, $false, $true, 1, $false, "69-25=", 2) | Out-Null
$d.Content.Find.Execute("11+40=", $true, $false, $false, $false, $false, $true, 1, $false, "34+36=", 2) | Out-Null
$d.Content.Find.Execute("8+9=", $true, $false, $false, $false, $false, $true, 1, $false, "2+7=", 2) | Out-Null
$d.Content.Find.Execute("79+16=", $true, $false, $false, $false, $false, $true, 1, $false, "43+22=", 2) | Out-Null
$d.Content.Find.Execute("16+29=", $true, $false, $false, $false, $false, $true, 1, $false, "94-56=", 2) | Out-Null
$d.Content.Find.Execute("40+39=", $true, $false, $false, $false, $false, $true, 1, $false, "52-29=", 2) | Out-Null
$d.Content.Find.Execute("51+2=", $true, $false, $false, $false, $false, $true, 1, $false, "11+88=", 2) | Out-Null
$d.Content.Find.Execute("70-6=", $true, $false, $false, $false, $false, $true, 1, $false, "60-25=", 2) | Out-Null
$d.Content.Find.Execute("99-34=", $true, $false, $false, $false, $false, $true, 1, $false, "80+10=", 2) | Out-Null
$d.Content.Find.Execute("34+8=", $true, $false, $false, $false, $false, $true, 1, $false, "19+10=", 2) | Out-Null
$d.Content.Find.Execute("63+24=", $true, $false, $false, $false, $false, $true, 1, $false, "61+35=", 2) | Out-Null
$d.Content.Find.Execute("34+25=", $true, $false, $false, $false, $false, $true, 1, $false, "15+26=", 2) | Out-Null
$d.Content.Find.Execute("9+64=", $true, $false, $false, $false, $false, $true, 1, $false, "64+2=", 2) | Out-Null
$d.Content.Find.Execute("49-3=", $true, $false, $false, $false, $false, $true, 1, $false, "29-1=", 2) | Out-Null
$d.Content.Find.Execute("86-23=", $true, $false, $false, $false, $false, $true, 1, $false, "80-68=", 2) | Out-Null
$d.Content.Find.Execute("92-84=", $true, $false, $false, $false, $false, $true, 1, $false, "40-13=", 2) | Out-Null
$d.Content.Find.Execute("39+14=", $true, $false, $false, $false, $false, $true, 1, $false, "5+92=", 2) | Out-Null
$d.Content.Find.Execute("21+45=", $true, $false, $false, $false, $false, $true, 1, $false, "7+10=", 2) | Out-Null
$d.Content.Find.Execute("82-48=", $true, $false, $false, $false, $false, $true, 1, $false, "10-8=", 2) | Out-Null
$d.Content.Find.Execute("5+36=", $true, $false, $false, $false, $false, $true, 1, $false, "82+9=", 2) | Out-Null
$d.Content.Find.Execute("42-32=", $true, $false, $false, $false, $false, $true, 1, $false, "47-7=", 2) | Out-Null
$d.Content.Find.Execute("4+79=", $true, $false, $false, $false, $false, $true, 1, $false, "32-3=", 2) | Out-Null
$d.Content.Find.Execute("64-27=", $true, $false, $false, $false, $false, $true, 1, $false, "45-39=", 2) | Out-Null
$d.Content.Find.Execute("11+66=", $true, $false, $false, $false, $false, $true, 1, $false, "94-8=", 2) | Out-Null
$d.Content.Find.Execute("49+27=", $true, $false, $false, $false, $false, $true, 1, $false, "72-37=", 2) | Out-Null
$d.Content.Find.Execute("11+22=", $true, $false, $false, $false, $false, $true, 1, $false, "67-16=", 2) | Out-Null
$d.Content.Find.Execute("6+28=", $true, $false, $false, $false, $false, $true, 1, $false, "75-53=", 2) | Out-Null
$d.Content.Find.Execute("9+30=", $true, $false, $false, $false, $false, $true, 1, $false, "54+39=", 2) | Out-Null
$d.Content.Find.Execute("5+69=", $true, $false, $false, $false, $false, $true, 1, $false, "47-27=", 2) | Out-Null
$d.Content.Find.Execute("2+83=", $true, $false, $false, $false, $false, $true, 1, $false, "8+28=", 2) | Out-Null
$d.Content.Find.Execute("72-14=", $true, $false, $false, $false, $false, $true, 1, $false, "39-5=", 2) | Out-Null
$d.Content.Find.Execute("6+51=", $true, $false, $false, $false, $false, $true, 1, $false, "3+55=", 2) | Out-Null
$d.Content.Find.Execute("48-14=", $true, $false, $false, $false, $false, $true, 1, $false, "44+55=", 2) | Out-Null
$d.Content.Find.Execute("5+2=", $true, $false, $false, $false, $false, $true, 1, $false, "35+10=", 2) | Out-Null
$d.Content.Find.Execute("28-5=", $true, $false, $false, $false, $false, $true, 1, $false, "29+28=", 2) | Out-Null
$d.Content.Find.Execute("13+23=", $true, $false, $false, $false, $false, $true, 1, $false, "92-60=", 2) | Out-Null
$d.Content.Find.Execute("94-74=", $true, $false, $false, $false, $false, $true, 1, $false, "71+21=", 2) | Out-Null
$d.Content.Find.Execute("70-69=", $true, $false, $false, $false, $false, $true, 1, $false, "25+1=", 2) | Out-Null
$d.Content.Find.Execute("38-0=", $true, $false, $false, $false, $false, $true, 1, $false, "25-9=", 2) | Out-Null
$d.Content.Find.Execute("56-14=", $true, $false, $false, $false, $false, $true, 1, $false, "81-47=", 2) | Out-Null
$d.Content.Find.Execute("88-35=", $true, $false, $false, $false, $false, $true, 1, $false, "16+79=", 2) | Out-Null
$d.Content.Find.Execute("91-75=", $true, $false, $false, $false, $false, $true, 1, $false, "34-28=", 2) | Out-Null
$d.Content.Find.Execute("97-29=", $true, $false, $false, $false, $false, $true, 1, $false, "61-30=", 2) | Out-Null
$d.Content.Find.Execute("9+32=", $true, $false, $false, $false, $false, $true, 1, $false, "93-29=", 2) | Out-Null
$d.Content.Find.Execute("20+20=", $true, $false, $false, $false, $false, $true, 1, $false, "92-61=", 2) | Out-Null
$d.Content.Find.Execute("99-13=", $true, $false, $false, $false, $false, $true, 1, $false, "11+63=", 2) | Out-Null
$d.Content.Find.Execute("96-58=", $true, $false, $false, $false, $false, $true, 1, $false, "98-68=", 2) | Out-Null
$d.Content.Find.Execute("56-17=", $true, $false, $false, $false, $false, $true, 1, $false, "98-84=", 2) | Out-Null
$d.Content.Find.Execute("27+45=", $true, $false, $false, $false, $false, $true, 1, $false, "84+9=", 2) | Out-Null
$d.Content.Find.Execute("21+13=", $true, $false, $false, $false, $false, $true, 1, $false, "72+15=", 2) | Out-Null
$d.Content.Find.Execute("15+63=", $true, $false, $false, $false, $false, $true, 1, $false, "53-45=", 2) | Out-Null
$d.Content.Find.Execute("36-0=", $true, $false, $false, $false, $false, $true, 1, $false, "84-23=", 2) | Out-Null
$d.Content.Find.Execute("90-12=", $true, $false, $false, $false, $false, $true, 1, $false, "45-36=", 2) | Out-Null
$d.Content.Find.Execute("61-52=", $true, $false, $false, $false, $false, $true, 1, $false, "64-7=", 2) | Out-Null
$d.Content.Find.Execute("60+18=", $true, $false, $false, $false, $false, $true, 1, $false, "12+46=", 2) | Out-Null
$d.Content.Find.Execute("0+26=", $true, $false, $false, $false, $false, $true, 1, $false, "89-50=", 2) | Out-Null
$d.Content.Find.Execute("48+0=", $true, $false, $false, $false, $false, $true, 1, $false, "14+45=", 2) | Out-Null
$d.Content.Find.Execute("62+1=", $true, $false, $false, $false, $false, $true, 1, $false, "39+48=", 2) | Out-Null
$d.Content.Find.Execute("90-23=", $true, $false, $false, $false, $false, $true, 1, $false, "79+20=", 2) | Out-Null
$d.Content.Find.Execute("23+22=", $true, $false, $false, $false, $false, $true, 1, $false, "78-42=", 2) | Out-Null
$d.Content.Find.Execute("9+57=", $true, $false, $false, $false, $false, $true, 1, $false, "87-6=", 2) | Out-Null
$d.Content.Find.Execute("39+25=", $true, $false, $false, $false, $false, $true, 1, $false, "39-32=", 2) | Out-Null
$d.Content.Find.Execute("77-47=", $true, $false, $false, $false, $false, $true, 1, $false, "56+5=", 2) | Out-Null
$d.Content.Find.Execute("5+7=", $true, $false, $false, $false, $false, $true, 1, $false, "18+21=", 2) | Out-Null
$d.Content.Find.Execute("46+43=", $true, $false, $false, $false, $false, $true, 1, $false, "8+23=", 2) | Out-Null
$d.Content.Find.Execute("53+27=", $true, $false, $false, $false, $false, $true, 1, $false, "27-4=", 2) | Out-Null
$d.Content.Find.Execute("98-3=", $true, $false, $false, $false, $false, $true, 1, $false, "6+72=", 2) | Out-Null
$d.Content.Find.Execute("43+23=", $true, $false, $false, $false, $false, $true, 1, $false, "18+64=", 2) | Out-Null
